$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H19: append new text to the existing "3 Big Data Challange (13:00)" note ---
# Final text (3 segments):
#   1) "3 Big Data Challange\n"                         -> default formatting
#   2) "(13:00)\n"                                       -> red, Calibri (Body), 12pt
#   3) "Determine how sediment will be included->upload plots due at 11:59pm"
#                                                         -> black/theme text, Calibri (Body), 12pt
$nl = [char]10
$part1 = "3 Big Data Challange" + $nl
$part2 = "(13:00)" + $nl
$part3 = "Determine how sediment will be included->upload plots due at 11:59pm"
$fullH19 = $part1 + $part2 + $part3

$h19 = $ws.Range("H19")
$h19.Value = $fullH19

$run2Start = $part1.Length + 1
$run2 = $h19.Characters($run2Start, $part2.Length)
$run2.Font.Name = "Calibri (Body)"
$run2.Font.Size = 12
$run2.Font.Color = 255

$run3Start = $part1.Length + $part2.Length + 1
$run3 = $h19.Characters($run3Start, $part3.Length)
$run3.Font.Name = "Calibri (Body)"
$run3.Font.Size = 12
$run3.Font.ColorIndex = 1

# row 19 now needs two wrapped lines of extra text -> taller row
$ws.Rows.Item(19).RowHeight = 68

# --- New assignment-due notes for the following two Wednesday/Thursday rows ---
$ws.Range("H20").Value = "Complete MANCOVA assumption tests DUE at 9pm"
$ws.Range("H21").Value = "MANCOVA, regression, Shapiro-Wilks due at 11:59pm"

# --- View state: zoom + active selection ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("H22").Select()
